$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin/Link text cells: the symbol list shifted down one row starting at row 17
# (HotbitToken inserted ahead of LEO, bumping LEO/GateToken/BitpandaEcosystemToken/
# ProBitToken/ZBToken/CoinExToken/BitKan down a row each). Plain text, safe to set directly. ---
$ws.Range("B17").Value = "HotbitToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"

# --- Price / Volume(1h) cells: these look numeric but must stay stored as TEXT
# (matching the workbook's original inlineStr cells) instead of being auto-coerced
# into floating point numbers/percentages by Excel. Force each cell to text format
# first, write the literal string, then restore the default 'Normal' style afterwards
# so no stray number-format styling lingers on the cells. ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D51:E51").NumberFormat = "@"

$ws.Range("D2").Value = "330.59"
$ws.Range("D3").Value = "45.37"
$ws.Range("E3").Value = "2.60%"
$ws.Range("D4").Value = "5.583"
$ws.Range("E4").Value = "1.65%"
$ws.Range("D5").Value = "0.08340"
$ws.Range("E5").Value = "3.81%"
$ws.Range("D6").Value = "2.104"
$ws.Range("E6").Value = "4.67%"
$ws.Range("D7").Value = "0.9781"
$ws.Range("E7").Value = "2.67%"
$ws.Range("E8").Value = "0.68%"
$ws.Range("D9").Value = "0.1199"
$ws.Range("E9").Value = "4.76%"
$ws.Range("D10").Value = "0.1922"
$ws.Range("E10").Value = "1.14%"
$ws.Range("D11").Value = "10.31"
$ws.Range("E11").Value = "-3.57%"
$ws.Range("D12").Value = "0.09846"
$ws.Range("E12").Value = "-0.24%"
$ws.Range("D13").Value = "0.04626"
$ws.Range("E13").Value = "-4.12%"
$ws.Range("D14").Value = "0.1057"
$ws.Range("E14").Value = "-0.71%"
$ws.Range("D15").Value = "0.001278"
$ws.Range("E15").Value = "-0.26%"
$ws.Range("D16").Value = "0.005914"
$ws.Range("E16").Value = "0.57%"
$ws.Range("D17").Value = "0.004581"
$ws.Range("E17").Value = "5.03%"
$ws.Range("D18").Value = "3.376"
$ws.Range("E18").Value = "0.26%"
$ws.Range("D19").Value = "4.449"
$ws.Range("E19").Value = "1.05%"
$ws.Range("D20").Value = "0.3342"
$ws.Range("E20").Value = "-2.60%"
$ws.Range("D21").Value = "0.1392"
$ws.Range("E21").Value = "-0.43%"
$ws.Range("D22").Value = "0.2784"
$ws.Range("E22").Value = "11.30%"
$ws.Range("D23").Value = "0.04173"
$ws.Range("E23").Value = "2.46%"
$ws.Range("D24").Value = "0.001293"
$ws.Range("E24").Value = "1.64%"
$ws.Range("D25").Value = "0.0001302"
$ws.Range("E25").Value = "8.56%"
$ws.Range("D26").Value = "0.0003747"
$ws.Range("E26").Value = "0.15%"
$ws.Range("D38").Value = "0.02703"
$ws.Range("E38").Value = "3.83%"
$ws.Range("D39").Value = "0.05745"
$ws.Range("E39").Value = "-0.64%"
$ws.Range("D40").Value = "0.007907"
$ws.Range("E40").Value = "4.33%"
$ws.Range("E41").Value = "1.80%"
$ws.Range("E42").Value = "5.25%"
$ws.Range("D43").Value = "0.002100"
$ws.Range("E43").Value = "4.25%"
$ws.Range("D44").Value = "0.008515"
$ws.Range("E44").Value = "-3.60%"
$ws.Range("D45").Value = "0.3371"
$ws.Range("D46").Value = "0.00007125"
$ws.Range("E46").Value = "2.13%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "0.27%"
$ws.Range("D49").Value = "0.003529"
$ws.Range("E49").Value = "0.94%"
$ws.Range("D50").Value = "0.003536"
$ws.Range("E50").Value = "0.27%"
$ws.Range("D51").Value = "0.00002104"
$ws.Range("E51").Value = "0.27%"

$ws.Range("D2").Style = "Normal"
$ws.Range("D3:E3").Style = "Normal"
$ws.Range("D4:E4").Style = "Normal"
$ws.Range("D5:E5").Style = "Normal"
$ws.Range("D6:E6").Style = "Normal"
$ws.Range("D7:E7").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9:E9").Style = "Normal"
$ws.Range("D10:E10").Style = "Normal"
$ws.Range("D11:E11").Style = "Normal"
$ws.Range("D12:E12").Style = "Normal"
$ws.Range("D13:E13").Style = "Normal"
$ws.Range("D14:E14").Style = "Normal"
$ws.Range("D15:E15").Style = "Normal"
$ws.Range("D16:E16").Style = "Normal"
$ws.Range("D17:E17").Style = "Normal"
$ws.Range("D18:E18").Style = "Normal"
$ws.Range("D19:E19").Style = "Normal"
$ws.Range("D20:E20").Style = "Normal"
$ws.Range("D21:E21").Style = "Normal"
$ws.Range("D22:E22").Style = "Normal"
$ws.Range("D23:E23").Style = "Normal"
$ws.Range("D24:E24").Style = "Normal"
$ws.Range("D25:E25").Style = "Normal"
$ws.Range("D26:E26").Style = "Normal"
$ws.Range("D38:E38").Style = "Normal"
$ws.Range("D39:E39").Style = "Normal"
$ws.Range("D40:E40").Style = "Normal"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43:E43").Style = "Normal"
$ws.Range("D44:E44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46:E46").Style = "Normal"
$ws.Range("D47:E47").Style = "Normal"
$ws.Range("D49:E49").Style = "Normal"
$ws.Range("D50:E50").Style = "Normal"
$ws.Range("D51:E51").Style = "Normal"
